$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Business Analyst(BPT) ---
$ws.Range("B2").Value = "Business Analyst(BPT)"
$ws.Range("D2").Value = "We are looking for a dynamic Business Analyst with 5–7 years of experience who can operate at an intermediate to senior level.`nLead requirements gathering`nDrive business process analysis and consulting efforts to deliver effective solutions`nTranslate business needs into clear and actionable user stories and documentation`n"
$ws.Range("D2").WrapText = $true

# --- Row 3: QA Analyst (existing row, new description) ---
$ws.Range("D3").Value = "Experience: 3–6 Years | Intermediate Level`nPerform functional and regression testing across various ServiceNow modules`nEnsure solution quality through effective test planning and execution"
$ws.Range("D3").WrapText = $true

# --- Row 4: new job - Infra Engineer(On-prem setup) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Infra Engineer(On-prem setup)"
$ws.Range("C4").Value = "Remote"
$ws.Range("D4").Value = "Responsible for the installation, configuration, and ongoing maintenance of self-hosted ServiceNow instances and MID Servers`nEnsure optimal performance, availability, and reliability of the platform infrastructure`nManage upgrades, patching, and environment tuning for enterprise-scale ServiceNow deployments"
$ws.Range("D4").WrapText = $true

# --- Row heights (wrapped multi-line descriptions) ---
$ws.Rows.Item(2).RowHeight = 72
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 57.6

# --- Selection moves to B4 (matches author's last click before upload) ---
$ws.Range("B4").Select()
